# Weekly update: insert a new price record at the top of the data block
# (row 159), pushing all existing records down by one row. The table grows
# from A1:R216 to A1:R217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 159; rows 159-216 shift down to 160-217.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row with this week's new record.
$ws.Cells.Item(159, 1).Value = 8
$ws.Cells.Item(159, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(159, 3).Value = "Coquimbo"
$ws.Cells.Item(159, 4).Value = 44588
$ws.Cells.Item(159, 5).Value = 4
$ws.Cells.Item(159, 6).Value = 100112012
$ws.Cells.Item(159, 7).Value = "Espinaca"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 2960
$ws.Cells.Item(159, 11).Value = 400
$ws.Cells.Item(159, 12).Value = 500
$ws.Cells.Item(159, 13).Value = 450
$ws.Cells.Item(159, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(159, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(159, 16).Value = 900
$ws.Cells.Item(159, 17).Value = 0.5
$ws.Cells.Item(159, 18).Value = "Hortaliza"
